$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 1193, shifting existing rows 1193:1258 down to 1194:1259
$ws.Rows.Item(1193).Insert()

# Populate the new row 1193 with the new record
$ws.Range("A1193").Value = 5
$ws.Range("B1193").Value = "Macroferia Regional de Talca"
$ws.Range("C1193").Value = "Maule"
$ws.Range("D1193").Value = 45267
$ws.Range("E1193").Value = 7
$ws.Range("F1193").Value = "Fruta"
$ws.Range("G1193").Value = 100104
$ws.Range("H1193").Value = "Frutos de pepita"
$ws.Range("I1193").Value = 100104005
$ws.Range("J1193").Value = "Pera"
$ws.Range("K1193").Value = "Packham's Triumph"
$ws.Range("L1193").Value = "Especial"
$ws.Range("M1193").Value = 280
$ws.Range("N1193").Value = 16000
$ws.Range("O1193").Value = 16000
$ws.Range("P1193").Value = 16000
$ws.Range("Q1193").Value = "$/bandeja 18 kilos granel"
$ws.Range("R1193").Value = "Provincia de Curicó"
$ws.Range("S1193").Value = 889
$ws.Range("T1193").Value = 18
